$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The report template gained two new parameter columns ("IS_Service_type",
# "IS_Product_type"). They go in between the existing "Portfolio" (B) and
# "Full_Project_name" (old C) columns, so insert two blank columns at C:D -
# this shifts "Full_Project_name" from C to E and the new cells inherit the
# bold/filled header style automatically.
$ws.Columns("C:D").Insert()

# Populate the two new header cells.
$ws.Range("C1").Value = "IS_Service_type"
$ws.Range("D1").Value = "IS_Product_type"

# Give the two new columns sensible widths (close to the template's).
$ws.Columns("C").ColumnWidth = 19.15
$ws.Columns("D").ColumnWidth = 22.65

# The existing AutoFilter only covered the original A1:C1 header; rebuild it
# over the full, now 5-column header range.
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:E1").AutoFilter()

# Keep the workbook-level hidden _FilterDatabase name (used by the
# AutoFilter) in sync with the new header range.
$fdName = $wb.Names.Item("_xlnm._FilterDatabase")
$fdName.RefersTo = '=' + $ws.Name + '!$A$1:$E$1'

# Leave the selection resting on the first data row, below the header.
$null = $ws.Range("A2").Select()

Write-Output "Inserted IS_Service_type / IS_Product_type parameter columns"
